# Insert two new rows of data at the top of the "Espinaca" price table
# (rows 209-210), shifting the existing rows 209-245 down to 211-247.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("209:210").Insert()

# New row 209: Primera quality, week of 2021-10-05
$ws.Cells.Item(209, 1).Value  = 9
$ws.Cells.Item(209, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(209, 3).Value  = "Metropolitana"
$ws.Cells.Item(209, 4).Value  = 44474
$ws.Cells.Item(209, 5).Value  = 13
$ws.Cells.Item(209, 6).Value  = 100112012
$ws.Cells.Item(209, 7).Value  = "Espinaca"
$ws.Cells.Item(209, 8).Value  = "Sin especificar"
$ws.Cells.Item(209, 9).Value  = "Primera"
$ws.Cells.Item(209, 10).Value = 250
$ws.Cells.Item(209, 11).Value = 6000
$ws.Cells.Item(209, 12).Value = 7000
$ws.Cells.Item(209, 13).Value = 6500
$ws.Cells.Item(209, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(209, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(209, 16).Value = 650
$ws.Cells.Item(209, 17).Value = 10
$ws.Cells.Item(209, 18).Value = "Hortaliza"

# New row 210: Segunda quality, week of 2021-10-05
$ws.Cells.Item(210, 1).Value  = 9
$ws.Cells.Item(210, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(210, 3).Value  = "Metropolitana"
$ws.Cells.Item(210, 4).Value  = 44474
$ws.Cells.Item(210, 5).Value  = 13
$ws.Cells.Item(210, 6).Value  = 100112012
$ws.Cells.Item(210, 7).Value  = "Espinaca"
$ws.Cells.Item(210, 8).Value  = "Sin especificar"
$ws.Cells.Item(210, 9).Value  = "Segunda"
$ws.Cells.Item(210, 10).Value = 106
$ws.Cells.Item(210, 11).Value = 4000
$ws.Cells.Item(210, 12).Value = 5000
$ws.Cells.Item(210, 13).Value = 4500
$ws.Cells.Item(210, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(210, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(210, 16).Value = 450
$ws.Cells.Item(210, 17).Value = 10
$ws.Cells.Item(210, 18).Value = "Hortaliza"
